$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Append new question rows (11-18) to the November tracker sheet.
# Column A = DATE (copies date format from existing A column),
# Column B = PROBLEM (copies fill style from existing B column),
# Column C = TOPICS (copies fill style from existing C column),
# Column D = RESOURCE (plain, no style).
# ---------------------------------------------------------------------------

$rows = @(
    @{ Row = 11; Date = 45964; Problem = "K Closest Pair to Origin";                   Topic = "HEAP";          Resource = "NEETCODE" },
    @{ Row = 12; Date = 45965; Problem = "Two Sum";                                    Topic = "ARRAY";         Resource = "NEETCODE" },
    @{ Row = 13; Date = 45965; Problem = "Find Maximum Ones";                          Topic = "ARRAY";         Resource = "NEETCODE" },
    @{ Row = 14; Date = 45966; Problem = "Minimum Number Games";                       Topic = "HEAP";          Resource = "LEETCODE" },
    @{ Row = 15; Date = 45966; Problem = "Connect Ropes To Minimize Cost";             Topic = "HEAP";          Resource = "INTERVIEW BIT" },
    @{ Row = 16; Date = 45967; Problem = "Binary Search for Descended Sorted Array";   Topic = "BINARY SEARCH"; Resource = "ADITYA VERMA" },
    @{ Row = 17; Date = 45967; Problem = "First and Last Occurrence of element";       Topic = "BINARY SEARCH"; Resource = "LEETCODE" },
    @{ Row = 18; Date = 45967; Problem = "Order Agnostic Search";                      Topic = "BINARY SEARCH"; Resource = "GFG" }
)

foreach ($r in $rows) {
    $rowIndex = $r.Row

    # Column A: date value, formatted like the rest of the DATE column.
    $ws.Range("A$rowIndex").Value = $r.Date
    $ws.Range("A2").Copy()
    $ws.Range("A$rowIndex").PasteSpecial(-4122)

    # Column B: problem name, formatted like the rest of the PROBLEM column.
    $ws.Range("B$rowIndex").Value = $r.Problem
    $ws.Range("B2").Copy()
    $ws.Range("B$rowIndex").PasteSpecial(-4122)

    # Column C: topic, formatted like the rest of the TOPICS column.
    $ws.Range("C$rowIndex").Value = $r.Topic
    $ws.Range("C2").Copy()
    $ws.Range("C$rowIndex").PasteSpecial(-4122)

    # Column D: resource, no special style.
    $ws.Range("D$rowIndex").Value = $r.Resource
}

# Row 17's PROBLEM cell also carries the DATE column's number format (merged
# with its existing fill), matching a stray formatting paste in the source.
$ws.Range("B17").NumberFormat = "m/d/yyyy"

# The PROBLEM column got noticeably wider to fit the new, longer questions.
$ws.Columns.Item(2).ColumnWidth = 43

# Selection ends up just past the newly-added data.
$ws.Range("A19").Select()
